$wb = $excel.ActiveWorkbook

# --- "Gen" sheet: append a new data row (row 10) ---------------------------
# Row 9 was the last data row (bus_i value 8); the new row continues the
# sequence with bus_i = 9 and zeros for every other column (B:U).
$genWs = $wb.Worksheets.Item("Gen")
$genWs.Range("A10").Value = 9
$genWs.Range("B10:U10").Value = 0

# --- Make "Gen" the active sheet/tab with a new selection -------------------
$genWs.Activate()
$genWs.Range("A11").Select()
